$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.448.64'
$ws.Range('E2').Value = '  -0.67%  '
$ws.Range('D3').Value = '2.202.40'
$ws.Range('E3').Value = '  -1.48%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '253.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.628'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '68.92'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.61%  '
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.586'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.64%  '
$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.13%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '36.76'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.90%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0952'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.10%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.79%  '
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').Value = '2.530.67'
$ws.Range('E15').Value = '  -1.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.896'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.97%  '
$ws.Range('D18').Value = '2.203.17'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').Value = '41.300.65'
$ws.Range('E19').Value = '  -0.72%  '
$ws.Range('D20').Value = '0.0₃0960'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.26'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -1.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.05'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +20.83%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.82'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.85%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('B29').Value = 'Monero'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.64'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.74%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.20%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.118'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.73'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.99%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0760'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.04%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.124'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.23'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.74%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.65'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.03'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0313'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.24'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.47'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.52%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '62.89'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.64%  '
$ws.Range('E44').Value = '  -4.44%  '
$ws.Range('E45').Value = '  -4.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.25'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +12.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.75'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.15%  '
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.01'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.44%  '
